$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.9     # Current Capital
$summary.Range("B4").Value = -0.11      # Total P&L $
$summary.Range("B6").Value = 98         # Total Trades
$summary.Range("B8").Value = 50         # Losing Trades
$summary.Range("B9").Value = 36.73      # Win Rate %

# --- Sheet: Strategy Status (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.90000000000001   # Capital
$status.Range("D4").Value = 98                   # Trades
$status.Range("E4").Value = -0.11                # P&L $
$status.Range("F4").Value = -0.1                 # P&L %
$status.Range("G4").Value = 36.73                # Win Rate %

# --- Sheet: All Trades - append new row 99 (Trade #98) ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A99").Value = 98
# Force text type so the date-like string isn't auto-converted to a date serial
$allTrades.Range("B99").NumberFormat = "@"
$allTrades.Range("B99").Value = "2026-02-17"
$allTrades.Range("B99").Style = "Normal"
$allTrades.Range("C99").Value = "15:57:27"
$allTrades.Range("D99").Value = "MarketMaking"
$allTrades.Range("E99").Value = "UP"
$allTrades.Range("F99").Value = 0.16
$allTrades.Range("G99").Value = 0.15
$allTrades.Range("H99").Value = "CLOSED"
$allTrades.Range("I99").Value = -6.25
$allTrades.Range("J99").Value = -0.01
$allTrades.Range("K99").Value = 99.90000000000001
$allTrades.Range("L99").Value = 0
$allTrades.Range("M99").Value = 0
$allTrades.Range("N99").Value = 0.6
$allTrades.Range("O99").Value = "Normal spread capture: 19600 bps"
$allTrades.Range("P99").Value = "early_exit"
$allTrades.Range("Q99").Value = 0.13

# --- Sheet: MarketMaking - append new row 99 (same trade) ---
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("A99").Value = 98
# Force text type so the date-like string isn't auto-converted to a date serial
$mm.Range("B99").NumberFormat = "@"
$mm.Range("B99").Value = "2026-02-17"
$mm.Range("B99").Style = "Normal"
$mm.Range("C99").Value = "15:57:27"
$mm.Range("D99").Value = "MarketMaking"
$mm.Range("E99").Value = "UP"
$mm.Range("F99").Value = 0.16
$mm.Range("G99").Value = 0.15
$mm.Range("H99").Value = "CLOSED"
$mm.Range("I99").Value = -6.25
$mm.Range("J99").Value = -0.01
$mm.Range("K99").Value = 99.90000000000001
$mm.Range("L99").Value = 0
$mm.Range("M99").Value = 0
$mm.Range("N99").Value = 0.6
$mm.Range("O99").Value = "Normal spread capture: 19600 bps"
$mm.Range("P99").Value = "early_exit"
$mm.Range("Q99").Value = 0.13
